{"js": "// Update the \"Port Hospital\" row in the transmission-chains table:\n//   Total cases:            1,365 -> 1,364\n//   Recovered / Total:        580 -> 579\n//   Recovered / % of cases:  42.5 -> 42.4\n//   Died / % of cases:       57.5 -> 57.6\n// Each value is unique in the document, so a simple search-and-replace\n// (matching the whole cell text, case-sensitive) is safe and unambiguous.\n\nconst replacements = [\n  { find: \"1,365\", replace: \"1,364\" },\n  { find: \"580\", replace: \"579\" },\n  { find: \"42.5\", replace: \"42.4\" },\n  { find: \"57.5\", replace: \"57.6\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: true });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${find}\" but found ${results.items.length}.`\n    );\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the \"Port Hospital\" row in the transmission-chains table:\n#   Total cases:            1,365 -> 1,364\n#   Recovered / Total:        580 -> 579\n#   Recovered / % of cases:  42.5 -> 42.4\n#   Died / % of cases:       57.5 -> 57.6\n# Each value is unique in the document, so Find/Replace against the whole\n# document content is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"1,365\"; Replace = \"1,364\" },\n    @{ Find = \"580\";    Replace = \"579\" },\n    @{ Find = \"42.5\";   Replace = \"42.4\" },\n    @{ Find = \"57.5\";   Replace = \"57.6\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $found = $find.Execute($r.Find, $false, $true, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        throw \"Could not find text '$($r.Find)' in document content.\"\n    }\n}\n"}
